$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.833.75"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "1.763.00"
$ws.Range("E3").Value = "  -2.55%  "

$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").Value = "321.68"
$ws.Range("E5").Value = "  -2.41%  "

$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").Value = "0.4254"
$ws.Range("E7").Value = "  -4.10%  "

$ws.Range("D8").Value = "0.3637"

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07491"
$ws.Range("E9").Value = "  -2.71%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "42.43"
$ws.Range("E10").Value = "  -5.14%  "

$ws.Range("D11").Value = "1.089"
$ws.Range("E11").Value = "  -2.61%  "

$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("E13").Value = "  -5.63%  "

$ws.Range("D14").Value = "6.067"
$ws.Range("E14").Value = "  -3.57%  "

$ws.Range("D15").Value = "7.290"
$ws.Range("E15").Value = "  -2.26%  "

$ws.Range("D16").Value = "1.780.34"
$ws.Range("E16").Value = "  -2.08%  "

$ws.Range("D17").Value = "91.23"
$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("E18").Value = "  -2.42%  "

$ws.Range("D19").Value = "0.06369"
$ws.Range("E19").Value = "  -1.69%  "

$ws.Range("D20").Value = "0.9989"

$ws.Range("D21").Value = "17.02"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").Value = "5.920"
$ws.Range("E22").Value = "  -5.40%  "

$ws.Range("D23").Value = "27.857.76"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  -3.93%  "

$ws.Range("D25").Value = "2.108"
$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("D26").Value = "157.34"
$ws.Range("E26").Value = "  +1.18%  "

$ws.Range("D27").Value = "20.23"
$ws.Range("E27").Value = "  -1.54%  "

$ws.Range("D28").Value = "1.968.58"
$ws.Range("E28").Value = "  -2.57%  "

$ws.Range("D29").Value = "2.135"
$ws.Range("E29").Value = "  -8.00%  "

$ws.Range("D30").Value = "124.86"
$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("D31").Value = "1.113"
$ws.Range("E31").Value = "  -7.33%  "

$ws.Range("D32").Value = "3.685"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").Value = "5.555"
$ws.Range("E33").Value = "  -5.01%  "

$ws.Range("D34").Value = "0.08869"
$ws.Range("E34").Value = "  -3.90%  "

$ws.Range("D35").Value = "12.23"
$ws.Range("E35").Value = "  -6.40%  "

$ws.Range("E36").Value = "  -2.21%  "

$ws.Range("D37").Value = "0.2103"
$ws.Range("E37").Value = "  -3.05%  "

$ws.Range("D38").Value = "0.06040"
$ws.Range("E38").Value = "  -2.59%  "

$ws.Range("D39").Value = "4.965"
$ws.Range("E39").Value = "  -3.84%  "

$ws.Range("D40").Value = "0.6323"

$ws.Range("D41").Value = "1.175"
$ws.Range("E41").Value = "  -1.59%  "

$ws.Range("D42").Value = "0.9981"
$ws.Range("E42").Value = "  -0.38%  "

$ws.Range("D43").Value = "7.876"
$ws.Range("E43").Value = "  -2.60%  "

$ws.Range("D44").Value = "1.396"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").Value = "13.29"
$ws.Range("E45").Value = "  -4.30%  "

$ws.Range("E46").Value = "  -3.29%  "

$ws.Range("D47").Value = "3.681"
$ws.Range("E47").Value = "  -2.17%  "

$ws.Range("D48").Value = "1.985"
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("D49").Value = "122.88"
$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("D50").Value = "1.185"
$ws.Range("E50").Value = "  +3.02%  "

$ws.Range("D51").Value = "0.06832"
$ws.Range("E51").Value = "  -2.07%  "
